# CRM.xlsx update: rename Sheet1 -> Signin, add a new SigninCustom sheet,
# populate sign-in test data (EMAIL/PASSWORD/RESULT) used by the Selenium
# DataProvider, add a mailto hyperlink on the admin2 row and switch the
# active tab to the newly added sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet1 -> "Signin" --------------------------------------------------
$signin = $wb.Worksheets.Item(1)
$signin.Name = "Signin"

# Header row
$signin.Range("A1").Value = "EMAIL"
$signin.Range("B1").Value = "PASSWORD"

# Data rows (typed before the RESULT column so the shared-string table
# order matches: PASSWORD, EMAIL, admin@demo.com, riseDemo, admin2@demo.com)
$signin.Range("A2").Value = "admin@demo.com"
$signin.Range("B2").Value = "riseDemo"
$signin.Range("A3").Value = "admin2@demo.com"
$signin.Range("B3").Value = "riseDemo"

# RESULT column (registered after the data rows above)
$signin.Range("C1").Value = "RESULT"
$signin.Range("C2").Value = "'"
$signin.Range("C3").Value = "'"

# Widen column A to fit the email addresses, drop the old bestFit autosize
$signin.Columns.Item(1).ColumnWidth = 29.5

# Hyperlink the second admin row's e-mail address
[void]$signin.Hyperlinks.Add($signin.Range("A3"), "mailto:admin2@demo.com")
$signin.Range("A3").Style = "Normal"

$signin.PageSetup.Orientation = 1

[void]$signin.Range("C4").Select()

# ---- New "SigninCustom" sheet --------------------------------------------
$custom = $wb.Worksheets.Add([Type]::Missing, $signin)
$custom.Name = "SigninCustom"

$custom.Range("A1").Value = "EMAIL"
$custom.Range("B1").Value = "PASSWORD"
$custom.Range("C1").Value = "RESULT"

$custom.Range("A2").Value = "admin@demo.com"
$custom.Range("B2").Value = "riseDemo"
$custom.Range("A3").Value = "admin@demo.com"
$custom.Range("B3").Value = "riseDemo"
$custom.Range("A4").Value = "client@demo.com"
$custom.Range("B4").Value = "riseDemo"
$custom.Range("A5").Value = "admin@demo.com"
$custom.Range("B5").Value = "riseDemo"
$custom.Range("A6").Value = "admin@demo.com"
$custom.Range("B6").Value = "riseDemo"

$custom.Range("C2").Value = "'"
$custom.Range("C3").Value = "'"
$custom.Range("C4").Value = "'"
$custom.Range("C5").Value = "'"
$custom.Range("C6").Value = "'"
$custom.Range("C7").Value = "'"

$custom.Columns.Item(1).ColumnWidth = 29.5

$custom.PageSetup.Orientation = 1

[void]$custom.Range("A8").Select()
